$wb = $excel.ActiveWorkbook
$wsInput = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

$newProductName = "2585-RBI-EI-DB-DL-REC-NON-RNI-CTPD-DL-MD-TR-2-DATE-VAR-INST-MIS-1st"

# Update product name on both sheets so the shared string is replaced in
# place rather than leaving the old text as an orphaned shared string.
$wsInput.Range("B1").Value = $newProductName
$wsOutput.Range("B1").Value = $newProductName

# Update short name (B2) from numeric 2585 to text "258e"
$wsInput.Range("B2").Value = "258e"

# Remove stale test selection (was on B13), reset to B1
$wsInput.Range("B1").Select()
